$wb = $excel.ActiveWorkbook

$wsObjects = $wb.Worksheets.Item("Objects")
$wsImages  = $wb.Worksheets.Item("Images")
$wsAlbums  = $wb.Worksheets.Item("Albums")

# =====================================================================
# 1) New catalog rows for Obj127-130 ("Objects" sheet, rows 66-69).
#    Written in the exact order the original strings were first typed
#    so the rebuilt shared-string table lines up with the source file.
# =====================================================================
$wsObjects.Range("H66").Value2 = "alb000127"
$wsObjects.Range("B66").Value2 = "2015-08-09-obj000127-01.textile"
$wsObjects.Range("G67").Value2 = "pic000231"
$wsImages.Range("B172").Value2 = "2015-08-09-pic000231.textile"
$wsImages.Range("C172").Value2 = "S 23-0018"
$wsObjects.Range("H67").Value2 = "alb000128"
$wsObjects.Range("F66").Value2 = " Eternity’s White Flag - Before - / And God- at every Gate-"
$wsObjects.Range("B67").Value2 = "2015-08-09-obj000128-01.textile"
$wsObjects.Range("F67").Value2 = "Jacket Night at the ID450 Collective"
$wsObjects.Range("G68").Value2 = "pic000232"
$wsImages.Range("C173").Value2 = "S 28-0007"
$wsObjects.Range("H68").Value2 = "alb000129"
$wsObjects.Range("F68").Value2 = "Terrible Scrabble"
$wsObjects.Range("B68").Value2 = "2015-08-09-obj000129-01.textile"
$wsImages.Range("B174").Value2 = "2015-08-09-pic000233.textile"
$wsObjects.Range("H69").Value2 = "alb000130"
$wsImages.Range("C174").Value2 = "S 43-0034"
$wsObjects.Range("F69").Value2 = "Listening to Dionne (1)"
$wsObjects.Range("B69").Value2 = "2015-08-09-obj000130-01.textile"
$wsObjects.Range("G69").Value2 = "pic000233"

# ---- repeat references to the same album / picture / textile names ----
$wsImages.Range("I171").Value2 = "alb000127"
$wsAlbums.Range("B34").Value2 = "alb000127"
$wsImages.Range("A172").Value2 = "pic000231"
$wsImages.Range("B173").Value2 = "2015-08-09-pic000231.textile"
$wsImages.Range("I172").Value2 = "alb000128"
$wsAlbums.Range("B35").Value2 = "alb000128"
$wsImages.Range("A173").Value2 = "pic000232"
$wsImages.Range("I173").Value2 = "alb000129"
$wsAlbums.Range("B36").Value2 = "alb000129"
$wsImages.Range("I174").Value2 = "alb000130"
$wsAlbums.Range("B37").Value2 = "alb000130"
$wsImages.Range("A174").Value2 = "pic000233"

# ---- remaining columns that reuse pre-existing shared strings ----
$wsObjects.Range("C66").Value2 = "art"
$wsObjects.Range("D66").Value2 = "artworks"
$wsObjects.Range("E66").Value2 = "Collage/photo"
$wsObjects.Range("G66").Value2 = "pic000230"
$wsObjects.Range("C67").Value2 = "art"
$wsObjects.Range("D67").Value2 = "artworks"
$wsObjects.Range("E67").Value2 = "Collage/photo"
$wsObjects.Range("C68").Value2 = "art"
$wsObjects.Range("D68").Value2 = "artworks"
$wsObjects.Range("E68").Value2 = "Collage/photo"
$wsObjects.Range("C69").Value2 = "art"
$wsObjects.Range("D69").Value2 = "artworks"
$wsObjects.Range("E69").Value2 = "Collage/photo"

# =====================================================================
# 2) Re-apply the explicit black font used on the pasted-in rows
#    (matches style index "3" = black Calibri 12 in the source file).
# =====================================================================
$wsObjects.Range("B67:G67").Font.Color = 0
$wsObjects.Range("B68:E68").Font.Color = 0
$wsObjects.Range("G68").Font.Color = 0
$wsObjects.Range("B69:E69").Font.Color = 0
$wsObjects.Range("G69").Font.Color = 0

$wsImages.Range("A172:I172").Font.Color = 0
$wsImages.Range("A173:I173").Font.Color = 0
$wsImages.Range("A174").Font.Color = 0
$wsImages.Range("I174").Font.Color = 0

# =====================================================================
# 3) View state: select/scroll each sheet, ending on "Objects" so it is
#    the active tab when the workbook is saved (matches activeTab drop).
# =====================================================================
$wsImages.Range("A175").Select() | Out-Null
$wsAlbums.Range("B37").Select() | Out-Null
$wsObjects.Range("A65").Select() | Out-Null
$wsObjects.Activate() | Out-Null
